$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column before E. This shifts the old E/F/G (Expected
#    Clusters / Expected Centers / Error Message) to F/G/H and carries their
#    formatting + the SALIDAS merge (E1:G1 -> F1:H1) along automatically.
$ws.Columns("E:E").Insert()

# 2. Grow the ENTRADAS merge (C1:D1) to also cover the new column E.
$ws.Range("C1:D1").UnMerge()
$ws.Range("C1:E1").Merge()

# 3. The continuation cells of the SALIDAS merge (now G1:H1) inherited the
#    old "blank" style; make them match F1's style instead (same as before
#    the edit, where the continuation cells mirrored the header cell style).
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 4. New column header in E2: "Number of iterations", bold + thin border,
#    no center/top alignment (unlike the other header cells).
$c = $ws.Range("E2")
$c.Value = "Number of iterations"
$c.ClearFormats()
$c.Font.Bold = $true
$c.Borders.LineStyle = 1

# 5. Fill in the new column's data for every test-case row.
$ws.Range("E3:E12").Value = 10

# 6. Restore the column widths for the new / shifted columns (values chosen so
#    the engine's pixel-quantized ColumnWidth lands on the saved workbook's
#    target width). The new column E loses its custom width on insert, so it
#    needs to be restored too.
$ws.Columns("E:E").ColumnWidth = 20.7
$ws.Columns("F:F").ColumnWidth = 17.3
$ws.Columns("G:G").ColumnWidth = 63.5
$ws.Columns("H:H").ColumnWidth = 48.9

# 7. Update the view: no frozen/scrolled top-left cell anymore, selection on E13.
$ws.Range("E13").Select()
